$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, styled like the other header cells (bold/border/center via style index 1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$timestamps = @(
    "2021-10-05 10:51:02.614673",
    "2021-10-05 10:51:02.614684",
    "2021-10-05 10:51:02.614687",
    "2021-10-05 10:51:02.614690",
    "2021-10-05 10:51:02.614693",
    "2021-10-05 10:51:02.614695",
    "2021-10-05 10:51:02.614698",
    "2021-10-05 10:51:02.614700",
    "2021-10-05 10:51:02.614703",
    "2021-10-05 10:51:02.614706",
    "2021-10-05 10:51:02.614708",
    "2021-10-05 10:51:02.614711",
    "2021-10-05 10:51:02.614713",
    "2021-10-05 10:51:02.614716",
    "2021-10-05 10:51:02.614718",
    "2021-10-05 10:51:02.614721",
    "2021-10-05 10:51:02.614724",
    "2021-10-05 10:51:02.614726",
    "2021-10-05 10:51:02.614729",
    "2021-10-05 10:51:02.614732",
    "2021-10-05 10:51:02.614734",
    "2021-10-05 10:51:02.614737",
    "2021-10-05 10:51:02.614739",
    "2021-10-05 10:51:02.614742",
    "2021-10-05 10:51:02.614744",
    "2021-10-05 10:51:02.614747",
    "2021-10-05 10:51:02.614749",
    "2021-10-05 10:51:02.614752",
    "2021-10-05 10:51:02.614754",
    "2021-10-05 10:51:02.614757",
    "2021-10-05 10:51:02.614759",
    "2021-10-05 10:51:02.614762",
    "2021-10-05 10:51:02.614764",
    "2021-10-05 10:51:02.614767",
    "2021-10-05 10:51:02.614769",
    "2021-10-05 10:51:02.614772",
    "2021-10-05 10:51:02.614774",
    "2021-10-05 10:51:02.614777",
    "2021-10-05 10:51:02.614779",
    "2021-10-05 10:51:02.614781",
    "2021-10-05 10:51:02.614784",
    "2021-10-05 10:51:02.614787",
    "2021-10-05 10:51:02.614789",
    "2021-10-05 10:51:02.614792",
    "2021-10-05 10:51:02.614794",
    "2021-10-05 10:51:02.614797",
    "2021-10-05 10:51:02.614799",
    "2021-10-05 10:51:02.614801",
    "2021-10-05 10:51:02.614804",
    "2021-10-05 10:51:02.614807",
    "2021-10-05 10:51:02.614809",
    "2021-10-05 10:51:02.614811",
    "2021-10-05 10:51:02.614814",
    "2021-10-05 10:51:02.614817",
    "2021-10-05 10:51:02.614819",
    "2021-10-05 10:51:02.614822"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
